# This workbook holds one daily price record per row for "Betarraga"
# (beet) at "Macroferia Regional de Talca" across rows 2..406. The
# author added one more day's worth of records: a new most-recent day
# was inserted at the top of this block (row 279, immediately after the
# last row that was already "settled"/shared with other sheets) and,
# to keep the existing rows intact, every row from 279..406 shifted
# down by one, with the last original row (406) now appearing twice:
# once shifted into 406 and once more, verbatim, as the brand new last
# row 407.
#
# Concretely: new_row[n] (for n = 280..406) = old_row[n-1] for the
# columns that actually vary row-to-row (D, I, J, K, L, M, O, P -- the
# rest, A/B/C/E/F/G/H/N/Q/R, are constant across this whole block).
# new_row[279] reuses old_row[406]'s I/J/K/L/M/O/P but gets a brand new
# date (44839) in D. new_row[407] is an exact copy of old_row[406]
# (every column, including D).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstRow = 279
$lastRow = 406
$newLastRow = 407

# Columns that vary row-to-row within this block.
$varyCols = @("D", "I", "J", "K", "L", "M", "O", "P")
$colIndex = @{ "A"=1; "B"=2; "C"=3; "D"=4; "E"=5; "F"=6; "G"=7; "H"=8; "I"=9; "J"=10; "K"=11; "L"=12; "M"=13; "N"=14; "O"=15; "P"=16; "Q"=17; "R"=18 }

# 1) Snapshot the *original* values of the varying columns for every row
#    in the block (279..406) before we start overwriting anything.
$snapshot = @{}
for ($r = $firstRow; $r -le $lastRow; $r++) {
    $rowVals = @{}
    foreach ($c in $varyCols) {
        $rowVals[$c] = $ws.Cells.Item($r, $colIndex[$c]).Value()
    }
    $snapshot[$r] = $rowVals
}

# Also snapshot the *entire* original last row (406) -- every column --
# since it gets copied verbatim into the new row 407.
$fullCols = @("A","B","C","D","E","F","G","H","I","J","K","L","M","N","O","P","Q","R")
$lastRowFull = @{}
foreach ($c in $fullCols) {
    $lastRowFull[$c] = $ws.Cells.Item($lastRow, $colIndex[$c]).Value()
}

# 2) Shift rows 280..406 down from 279..405 (write highest row first isn't
#    required since we're writing from a snapshot, not live cells).
for ($r = $lastRow; $r -ge ($firstRow + 1); $r--) {
    $src = $snapshot[$r - 1]
    foreach ($c in $varyCols) {
        $ws.Cells.Item($r, $colIndex[$c]).Value = $src[$c]
    }
}

# 3) New row 279: I/J/K/L/M/O/P come from the old last row (406); D gets
#    a brand new date.
$srcLast = $snapshot[$lastRow]
foreach ($c in $varyCols) {
    if ($c -eq "D") {
        $ws.Cells.Item($firstRow, $colIndex[$c]).Value = 44839
    } else {
        $ws.Cells.Item($firstRow, $colIndex[$c]).Value = $srcLast[$c]
    }
}

# 4) Append brand new row 407, an exact copy of the original row 406.
foreach ($c in $fullCols) {
    $ws.Cells.Item($newLastRow, $colIndex[$c]).Value = $lastRowFull[$c]
}
$ws.Cells.Item($newLastRow, $colIndex["D"]).NumberFormat = $ws.Cells.Item($lastRow, $colIndex["D"]).NumberFormat()

Write-Host "done"
